$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.03904870566557
$ws.Range("D2").Value = 1.039573821909864
$ws.Range("E2").Value = 1.046505508868202
$ws.Range("F2").Value = 1.054732616899275
$ws.Range("I2").Value = 1.035355289966702
$ws.Range("J2").Value = 1.044142831731078
$ws.Range("K2").Value = 1.042358624097512
$ws.Range("L2").Value = 1.049270756065966
$ws.Range("M2").Value = 1.057475016126918
$ws.Range("N2").Value = 1.018599691489917

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.040262663587245
$ws.Range("D3").Value = 1.040638363080612
$ws.Range("E3").Value = 1.04762084730265
$ws.Range("F3").Value = 1.056021930560533
$ws.Range("I3").Value = 1.035582670209815
$ws.Range("J3").Value = 1.045000510791169
$ws.Range("K3").Value = 1.043232741507759
$ws.Range("L3").Value = 1.050196961774626
$ws.Range("M3").Value = 1.058576415766
$ws.Range("N3").Value = 1.018886367530547

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.041047314531491
$ws.Range("D4").Value = 1.041326690368995
$ws.Range("E4").Value = 1.04834219839105
$ws.Range("F4").Value = 1.0568561116763
$ws.Range("I4").Value = 1.035727358664905
$ws.Range("J4").Value = 1.045554125875561
$ws.Range("K4").Value = 1.043797258629948
$ws.Range("L4").Value = 1.050795346466586
$ws.Range("M4").Value = 1.05928845290075
$ws.Range("N4").Value = 1.019071333369091

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.041376978058296
$ws.Range("D5").Value = 1.041615944383286
$ws.Range("E5").Value = 1.048645372791096
$ws.Range("F5").Value = 1.057206781726231
$ws.Range("I5").Value = 1.035787601331858
$ws.Range("J5").Value = 1.045786541267573
$ws.Range("K5").Value = 1.044034320710817
$ws.Range("L5").Value = 1.051046685627428
$ws.Range("M5").Value = 1.05958764138645
$ws.Range("N5").Value = 1.019148965866174

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.041432318137519
$ws.Range("D6").Value = 1.041664504480542
$ws.Range("E6").Value = 1.048696272367747
$ws.Range("F6").Value = 1.057265659676778
$ws.Range("I6").Value = 1.03579768207925
$ws.Range("J6").Value = 1.045825545850985
$ws.Range("K6").Value = 1.044074109214416
$ws.Range("L6").Value = 1.051088873615689
$ws.Range("M6").Value = 1.059637867592404
$ws.Range("N6").Value = 1.019161993259774

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.041051720309236
$ws.Range("D7").Value = 1.041330555860283
$ws.Range("E7").Value = 1.048346249740899
$ws.Range("F7").Value = 1.056860797421448
$ws.Range("I7").Value = 1.035728165926282
$ws.Range("J7").Value = 1.045557232694396
$ws.Range("K7").Value = 1.043800427288721
$ws.Range("L7").Value = 1.050798705742374
$ws.Range("M7").Value = 1.059292451266706
$ws.Range("N7").Value = 1.019072371197415

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.039459148153904
$ws.Range("D8").Value = 1.039933693385955
$ws.Range("E8").Value = 1.046882515616856
$ws.Range("F8").Value = 1.055168366608406
$ws.Range("I8").Value = 1.035432639896037
$ws.Range("J8").Value = 1.044432970925531
$ws.Range("K8").Value = 1.042654263743004
$ws.Range("L8").Value = 1.049583965539049
$ws.Range("M8").Value = 1.05784737332081
$ws.Range("N8").Value = 1.018696685544703

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.036646129853127
$ws.Range("D9").Value = 1.037468333922469
$ws.Range("E9").Value = 1.044300495865185
$ws.Range("F9").Value = 1.052185294810353
$ws.Range("I9").Value = 1.034893174035429
$ws.Range("J9").Value = 1.042441395888919
$ws.Range("K9").Value = 1.040626129234296
$ws.Range("L9").Value = 1.047436232003699
$ws.Range("M9").Value = 1.055295962564934
$ws.Range("N9").Value = 1.018030581259352

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.034766122767598
$ws.Range("D10").Value = 1.035822038917788
$ws.Range("E10").Value = 1.042577207388241
$ws.Range("F10").Value = 1.050195910798633
$ws.Range("I10").Value = 1.034520932776569
$ws.Range("J10").Value = 1.041106542230809
$ws.Range("K10").Value = 1.039268273446667
$ws.Range("L10").Value = 1.045999469802657
$ws.Range("M10").Value = 1.053591544179512
$ws.Range("N10").Value = 1.017583729843039

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.033950917731294
$ws.Range("D11").Value = 1.035108507426771
$ws.Range("E11").Value = 1.041830519555309
$ws.Range("F11").Value = 1.049334295290027
$ws.Range("I11").Value = 1.034356755672148
$ws.Range("J11").Value = 1.040526822130337
$ws.Range("K11").Value = 1.038678919042832
$ws.Range("L11").Value = 1.045376142135908
$ws.Range("M11").Value = 1.052852660220161
$ws.Range("N11").Value = 1.017389572144134

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.033647938387484
$ws.Range("D12").Value = 1.034843366652748
$ws.Range("E12").Value = 1.041553090237034
$ws.Range("F12").Value = 1.049014220183822
$ws.Range("I12").Value = 1.03429532264083
$ws.Range("J12").Value = 1.040311227986485
$ws.Range("K12").Value = 1.038459795300717
$ws.Range("L12").Value = 1.045144428157694
$ws.Range("M12").Value = 1.052578074229238
$ws.Range("N12").Value = 1.017317352372196

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.033712936496341
$ws.Range("D13").Value = 1.034900244960823
$ws.Range("E13").Value = 1.041612603251982
$ws.Range("F13").Value = 1.049082878907983
$ws.Range("I13").Value = 1.034308520621866
$ws.Range("J13").Value = 1.040357485471878
$ws.Range("K13").Value = 1.038506807677973
$ws.Range("L13").Value = 1.045194139876011
$ws.Range("M13").Value = 1.052636979850697
$ws.Range("N13").Value = 1.017332848336449

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.033925876970976
$ws.Range("D14").Value = 1.035086592928902
$ws.Range("E14").Value = 1.041807588728856
$ws.Range("F14").Value = 1.049307838466293
$ws.Range("I14").Value = 1.034351686788474
$ws.Range("J14").Value = 1.040509006372142
$ws.Range("K14").Value = 1.038660810526741
$ws.Range("L14").Value = 1.045356992317398
$ws.Range("M14").Value = 1.05282996555514
$ws.Range("N14").Value = 1.017383604498172

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.03405705322158
$ws.Range("D15").Value = 1.035201394299087
$ws.Range("E15").Value = 1.041927715551323
$ws.Range("F15").Value = 1.049446439050988
$ws.Range("I15").Value = 1.034378223198037
$ws.Range("J15").Value = 1.040602328851107
$ws.Range("K15").Value = 1.038755668712563
$ws.Range("L15").Value = 1.045457306860209
$ws.Range("M15").Value = 1.052948852906371
$ws.Range("N15").Value = 1.017414863649077

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.034820200749569
$ws.Range("D16").Value = 1.035869379262347
$ws.Range("E16").Value = 1.042626752034868
$ws.Range("F16").Value = 1.050253088946362
$ws.Range("I16").Value = 1.034531765543739
$ws.Range("J16").Value = 1.041144979949379
$ws.Range("K16").Value = 1.039307357431583
$ws.Range("L16").Value = 1.046040812549239
$ws.Range("M16").Value = 1.0536405631559
$ws.Range("N16").Value = 1.017596601324456

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.03529859301808
$ws.Range("D17").Value = 1.036288206479314
$ws.Range("E17").Value = 1.043065105728307
$ws.Range("F17").Value = 1.050759023817251
$ws.Range("I17").Value = 1.034627276622619
$ws.Range("J17").Value = 1.041484908891035
$ws.Range("K17").Value = 1.03965304269942
$ws.Range("L17").Value = 1.04640650751689
$ws.Range("M17").Value = 1.054074222869731
$ws.Range("N17").Value = 1.017710421289804

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.035577520436321
$ws.Range("D18").Value = 1.036532436330292
$ws.Range("E18").Value = 1.043320742649048
$ws.Range("F18").Value = 1.051054108050723
$ws.Range("I18").Value = 1.03468269775924
$ws.Range("J18").Value = 1.041683017901956
$ws.Range("K18").Value = 1.039854540705514
$ws.Range("L18").Value = 1.046619695545731
$ws.Range("M18").Value = 1.054327086477617
$ws.Range("N18").Value = 1.017776746129503

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.03567260879704
$ws.Range("D19").Value = 1.036615701436602
$ws.Range("E19").Value = 1.043407900263683
$ws.Range("F19").Value = 1.051154721007822
$ws.Range("I19").Value = 1.034701545953705
$ws.Range("J19").Value = 1.041750539891842
$ws.Range("K19").Value = 1.039923223583796
$ws.Range("L19").Value = 1.046692367596937
$ws.Range("M19").Value = 1.054413292468105
$ws.Range("N19").Value = 1.017799350276481

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.035247277529503
$ws.Range("D20").Value = 1.036243277019332
$ws.Range("E20").Value = 1.043018079442552
$ws.Range("F20").Value = 1.050704743782459
$ws.Range("I20").Value = 1.03461705906994
$ws.Range("J20").Value = 1.041448454895936
$ws.Range("K20").Value = 1.039615967853801
$ws.Range("L20").Value = 1.046367283886471
$ws.Range("M20").Value = 1.05402770386631
$ws.Range("N20").Value = 1.01769821615256

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.033863176189918
$ws.Range("D21").Value = 1.035031720967756
$ws.Range("E21").Value = 1.041750172487452
$ws.Range("F21").Value = 1.04924159441004
$ws.Range("I21").Value = 1.034338987864126
$ws.Range("J21").Value = 1.040464394444338
$ws.Range("K21").Value = 1.038615466373473
$ws.Range("L21").Value = 1.045309041381894
$ws.Range("M21").Value = 1.052773139727537
$ws.Range("N21").Value = 1.017368660865268

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.032991917465093
$ws.Range("D22").Value = 1.034269368290551
$ws.Range("E22").Value = 1.040952547019436
$ws.Range("F22").Value = 1.048321463469906
$ws.Range("I22").Value = 1.034161547644745
$ws.Range("J22").Value = 1.039844169674047
$ws.Range("K22").Value = 1.037985188072436
$ws.Range("L22").Value = 1.044642626403314
$ws.Range("M22").Value = 1.051983582893688
$ws.Range("N22").Value = 1.017160871971382

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.033453885886839
$ws.Range("D23").Value = 1.034673563354111
$ws.Range("E23").Value = 1.041375426016115
$ws.Range("F23").Value = 1.048809261332943
$ws.Range("I23").Value = 1.034255859235471
$ws.Range("J23").Value = 1.040173106046655
$ws.Range("K23").Value = 1.038319427107384
$ws.Range("L23").Value = 1.044996006360851
$ws.Range("M23").Value = 1.052402215017313
$ws.Range("N23").Value = 1.017271080387632

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.035270465131909
$ws.Range("D24").Value = 1.036263578909132
$ws.Range("E24").Value = 1.043039328743248
$ws.Range("F24").Value = 1.050729270652462
$ws.Range("I24").Value = 1.034621676834751
$ws.Range("J24").Value = 1.041464927399438
$ws.Range("K24").Value = 1.039632720795832
$ws.Range("L24").Value = 1.046385007714061
$ws.Range("M24").Value = 1.054048724057287
$ws.Range("N24").Value = 1.017703731327774

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.037374172693537
$ws.Range("D25").Value = 1.038106161274087
$ws.Range("E25").Value = 1.044968344427688
$ws.Range("F25").Value = 1.052956598280683
$ws.Range("I25").Value = 1.035034857957015
$ws.Range("J25").Value = 1.042957517085732
$ws.Range("K25").Value = 1.041151460542399
$ws.Range("L25").Value = 1.047992336259968
$ws.Range("M25").Value = 1.055956167753448
$ws.Range("N25").Value = 1.018203273542583
